$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("B2").Value = "Ingeniería de Sistemas"
$ws.Range("B3").Value = "Ingeniería de Software"
$ws.Range("B4").Value = "Ciencias de la Computación"

$ws.Range("B5").Select()
